# Fruta / hortaliza, semanal
# A new daily price record is inserted as row 902 (pushing the existing
# rows 902..1002 down to 903..1003, and extending the used range to T1003).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 902; Excel shifts rows 902..1002 down to 903..1003
$ws.Rows.Item(902).Insert()

# Populate the newly inserted row 902 with the new record's data
$ws.Range('A902').Value = 9
$ws.Range('B902').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C902').Value = 'Metropolitana'
$ws.Range('D902').Value = 44918
$ws.Range('E902').Value = 13
$ws.Range('F902').Value = 'Fruta'
$ws.Range('G902').Value = 100104
$ws.Range('H902').Value = 'Frutos de pepita'
$ws.Range('I902').Value = 100104005
$ws.Range('J902').Value = 'Pera'
$ws.Range('K902').Value = "Packham's Triumph"
$ws.Range('L902').Value = 'Primera'
$ws.Range('M902').Value = 220
$ws.Range('N902').Value = 23000
$ws.Range('O902').Value = 23000
$ws.Range('P902').Value = 23000
$ws.Range('Q902').Value = '$/caja 18 kilos granel'
$ws.Range('R902').Value = 'Provincia de Curicó'
$ws.Range('S902').Value = 1278
$ws.Range('T902').Value = 18
